$d = $word.ActiveDocument

# Locate the anchor paragraph: "This keeps repeating each time a user clicks on an answer"
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "This keeps repeating each time a user clicks on an answer*") {
        $anchorIndex = $i
        break
    }
}

# Insert first new sub-bullet (ilvl 1) before the anchor paragraph.
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$p1 = $d.Paragraphs($anchorIndex)
$p1.Range.Text = "Depending on which answer is clicked, a variable for one of the four houses is increased by one"
$p1.Range.ListFormat.ListLevelNumber = 2
$anchorIndex = $anchorIndex + 1

# Insert second new sub-bullet (ilvl 1) before the anchor paragraph.
$d.Paragraphs($anchorIndex).Range.InsertParagraphBefore()
$p2 = $d.Paragraphs($anchorIndex)
$p2.Range.Text = "This means I need to somehow link a given answer option to a certain house variable"
$p2.Range.ListFormat.ListLevelNumber = 2
$anchorIndex = $anchorIndex + 1

# $anchorIndex now points back at "This keeps repeating...". Insert a new
# sub-bullet (ilvl 1) directly after it.
$d.Paragraphs($anchorIndex).Range.InsertParagraphAfter()
$p3 = $d.Paragraphs($anchorIndex + 1)
$p3.Range.Text = "This means that each question and its four “answer choices” needs to be associated with a state variable – each time the state variable value changes, so does the question that is displayed."
$p3.Range.ListFormat.ListLevelNumber = 2
